$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.924.16"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "3.128.24"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "588.29"
$ws.Range("E5").Value = "  -2.49%  "

$ws.Range("D6").Value = "135.84"
$ws.Range("E6").Value = "  -4.85%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "3.122.47"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("E10").Value = "  -4.01%  "

$ws.Range("D11").Value = "5.22"
$ws.Range("E11").Value = "  -2.97%  "

$ws.Range("E12").Value = "  -3.31%  "

$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  -5.56%  "

$ws.Range("D14").Value = "33.89"
$ws.Range("E14").Value = "  -3.59%  "

$ws.Range("D15").Value = "3.635.51"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").Value = "62.999.56"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("D18").Value = "3.118.82"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  -4.23%  "

$ws.Range("D20").Value = "469.35"
$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  -3.60%  "

$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("E23").Value = "  -0.52%  "

$ws.Range("D24").Value = "85.17"

$ws.Range("D25").Value = "12.90"
$ws.Range("E25").Value = "  -4.06%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("E28").Value = "  -5.94%  "

$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("D30").Value = "6.80"
$ws.Range("E30").Value = "  -4.91%  "

$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").Value = "26.55"
$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("E33").Value = "  -4.43%  "

$ws.Range("E34").Value = "  -4.42%  "

$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("D36").Value = "51.96"
$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("D37").Value = "5.72"
$ws.Range("E37").Value = "  -4.16%  "

$ws.Range("D38").Value = "0.0₃0675"
$ws.Range("E38").Value = "  -12.43%  "

$ws.Range("E39").Value = "  -2.24%  "

$ws.Range("D40").Value = "416.21"
$ws.Range("E40").Value = "  -6.39%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.905.06"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "8.15"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("E43").Value = "  -11.50%  "

$ws.Range("E44").Value = "  -6.79%  "

$ws.Range("D45").Value = "0.257"
$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  -6.09%  "

$ws.Range("D48").Value = "25.22"
$ws.Range("E48").Value = "  -3.02%  "

$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  -8.30%  "

$ws.Range("D51").Value = "120.51"
$ws.Range("E51").Value = "  -0.01%  "
